# Day 4 PPT.pptx edit script
# 1. Move the "Units of Measurement" slide (SlideID 259) from position 6 to position 3
# 2. Fix a typo on the "Box Model" content slide: "groupor" -> "group or"

$p = $ppt.ActivePresentation

# --- 1. Reorder slides -------------------------------------------------
# Find the slide whose SlideID is 259 ("Units of Measurement") and move it
# so that it becomes the 3rd slide in the deck.
$targetSlideId = 259
$destIndex = 3

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq $targetSlideId) {
        if ($i -ne $destIndex) {
            $slide.MoveTo($destIndex)
        }
        break
    }
}

# --- 2. Fix typo on the "Box Model" slide -------------------------------
# Locate the slide that contains the "groupor" typo and correct it to
# "group or".
foreach ($slide in $p.Slides) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*groupor*") {
                $tr.Replace("groupor", "group or") | Out-Null
            }
        }
    }
}
